$wb = $excel.ActiveWorkbook

# Rename second sheet from "Sheet1" to "R1R2"
$sheetR1R2 = $wb.Worksheets.Item("Sheet1")
$sheetR1R2.Name = "R1R2"

# Update R1/R2 resistor values on the R1R2 sheet
$sheetR1R2.Range("B1").Value = 205
$sheetR1R2.Range("B2").Value = 54.9

# Adjust the R1R2 sheet selection
$sheetR1R2.Activate() | Out-Null
$sheetR1R2.Range("E10").Select() | Out-Null

# Adjust the main sheet view (zoom + selection)
$mainSheet = $wb.Worksheets.Item("TPS62147")
$mainSheet.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 55
$mainSheet.Range("A1:A8").Select() | Out-Null
